$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) column writes to stay text-typed (matches the
# original inlineStr cells) instead of being auto-parsed as numbers,
# while preserving each cells original style/format.

$s_D2 = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.399.62"
$ws.Range("D2").Style = $s_D2
$ws.Range("E2").Value = "  +0.04%  "

$s_D3 = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.572.10"
$ws.Range("D3").Style = $s_D3
$ws.Range("E3").Value = "  +0.21%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("E5").Value = "  +0.07%  "

$s_D6 = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "290.31"
$ws.Range("D6").Style = $s_D6
$ws.Range("E6").Value = "  -0.36%  "

$s_D7 = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3765"
$ws.Range("D7").Style = $s_D7
$ws.Range("E7").Value = "  +2.75%  "

$s_D8 = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.86"
$ws.Range("D8").Style = $s_D8
$ws.Range("E8").Value = "  +0.92%  "

$s_D9 = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3426"
$ws.Range("D9").Style = $s_D9
$ws.Range("E9").Value = "  +1.27%  "

$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$s_D10 = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07649"
$ws.Range("D10").Style = $s_D10
$ws.Range("E10").Value = "  +0.78%  "

$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$s_D11 = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.157"
$ws.Range("D11").Style = $s_D11
$ws.Range("E11").Value = "  -1.35%  "

$s_D12 = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = $s_D12
$ws.Range("E12").Value = "  +0.10%  "

$s_D13 = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.25"
$ws.Range("D13").Style = $s_D13
$ws.Range("E13").Value = "  +0.28%  "

$s_D14 = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.027"
$ws.Range("D14").Style = $s_D14
$ws.Range("E14").Value = "  -0.52%  "

$s_D15 = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.943"
$ws.Range("D15").Style = $s_D15
$ws.Range("E15").Value = "  +0.69%  "

$s_D16 = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.572.01"
$ws.Range("D16").Style = $s_D16
$ws.Range("E16").Value = "  +0.13%  "

$s_D17 = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001134"
$ws.Range("D17").Style = $s_D17
$ws.Range("E17").Value = "  -0.53%  "

$s_D18 = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.32"
$ws.Range("D18").Style = $s_D18
$ws.Range("E18").Value = "  +1.42%  "

$s_D19 = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06758"
$ws.Range("D19").Style = $s_D19
$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("E20").Value = "  +0.04%  "

$ws.Range("E21").Value = "  +2.32%  "

$s_D22 = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.207"
$ws.Range("D22").Style = $s_D22
$ws.Range("E22").Value = "  -0.45%  "

$s_D23 = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.01"
$ws.Range("D23").Style = $s_D23
$ws.Range("E23").Value = "  -0.44%  "

$ws.Range("B24").Value = "WrappedBTC"
$ws.Range("C24").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$s_D24 = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.389.31"
$ws.Range("D24").Style = $s_D24
$ws.Range("E24").Value = "  -0.04%  "

$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$s_D25 = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.417"
$ws.Range("D25").Style = $s_D25
$ws.Range("E25").Value = "  +0.67%  "

$s_D26 = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.693"
$ws.Range("D26").Style = $s_D26
$ws.Range("E26").Value = "  -10.27%  "

$s_D27 = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.24"
$ws.Range("D27").Style = $s_D27
$ws.Range("E27").Value = "  +1.79%  "

$s_D28 = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "147.47"
$ws.Range("D28").Style = $s_D28
$ws.Range("E28").Value = "  +1.86%  "

$s_D29 = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.036"
$ws.Range("D29").Style = $s_D29
$ws.Range("E29").Value = "  +1.50%  "

$s_D30 = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.35"
$ws.Range("D30").Style = $s_D30
$ws.Range("E30").Value = "  +0.85%  "

$s_D31 = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.746.47"
$ws.Range("D31").Style = $s_D31
$ws.Range("E31").Value = "  +0.13%  "

$s_D32 = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.172"
$ws.Range("D32").Style = $s_D32
$ws.Range("E32").Value = "  -1.43%  "

$s_D33 = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.009"
$ws.Range("D33").Style = $s_D33
$ws.Range("E33").Value = "  +1.26%  "

$s_D34 = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9947"
$ws.Range("D34").Style = $s_D34
$ws.Range("E34").Value = "  -4.26%  "

$s_D35 = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.985"
$ws.Range("D35").Style = $s_D35
$ws.Range("E35").Value = "  -3.49%  "

$s_D36 = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08555"
$ws.Range("D36").Style = $s_D36
$ws.Range("E36").Value = "  +1.42%  "

$s_D37 = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02544"
$ws.Range("D37").Style = $s_D37
$ws.Range("E37").Value = "  -0.78%  "

$s_D38 = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2317"
$ws.Range("D38").Style = $s_D38
$ws.Range("E38").Value = "  +0.51%  "

$s_D39 = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06578"
$ws.Range("D39").Style = $s_D39
$ws.Range("E39").Value = "  +0.63%  "

$s_D40 = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.324"
$ws.Range("D40").Style = $s_D40
$ws.Range("E40").Value = "  +5.64%  "

$s_D41 = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.437"
$ws.Range("D41").Style = $s_D41
$ws.Range("E41").Value = "  -1.51%  "

$s_D42 = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.50"
$ws.Range("D42").Style = $s_D42

$s_D43 = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6418"
$ws.Range("D43").Style = $s_D43
$ws.Range("E43").Value = "  +0.43%  "

$ws.Range("E44").Value = "  +0.13%  "

$s_D45 = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.99"
$ws.Range("D45").Style = $s_D45
$ws.Range("E45").Value = "  -3.02%  "

$s_D46 = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.783"
$ws.Range("D46").Style = $s_D46
$ws.Range("E46").Value = "  +0.01%  "

$s_D47 = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5989"
$ws.Range("D47").Style = $s_D47
$ws.Range("E47").Value = "  -0.35%  "

$s_D48 = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.307"
$ws.Range("D48").Style = $s_D48
$ws.Range("E48").Value = "  +7.75%  "

$s_D49 = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.086"
$ws.Range("D49").Style = $s_D49
$ws.Range("E49").Value = "  -2.24%  "

$ws.Range("E50").Value = "  +2.03%  "

$s_D51 = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07324"
$ws.Range("D51").Style = $s_D51
$ws.Range("E51").Value = "  +0.48%  "
